$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: ref -> new value (numeric-looking text values are
# prefixed with a leading apostrophe so Excel stores them as literal
# text instead of auto-converting to numbers/percentages, matching
# the original inline-string cell contents.)
$updates = [ordered]@{
    'D2' = '''332.72'
    'E2' = '''1.91%'
    'G2' = '''23'
    'D3' = '''44.91'
    'E3' = '''2.21%'
    'G3' = '''23'
    'D4' = '''5.550'
    'E4' = '''-0.29%'
    'G4' = '''23'
    'D5' = '''0.08321'
    'E5' = '''3.52%'
    'G5' = '''23'
    'E6' = '''5.09%'
    'G6' = '''23'
    'D7' = '''0.9805'
    'E7' = '''3.87%'
    'G7' = '''23'
    'D8' = '''0.1124'
    'E8' = '''-3.76%'
    'G8' = '''23'
    'D9' = '''0.1915'
    'E9' = '''3.38%'
    'G9' = '''23'
    'D10' = '''10.33'
    'E10' = '''-14.70%'
    'G10' = '''23'
    'D11' = '''0.1006'
    'E11' = '''2.63%'
    'G11' = '''23'
    'D12' = '''0.04626'
    'E12' = '''-3.12%'
    'G12' = '''23'
    'D13' = '''0.1058'
    'E13' = '''-0.74%'
    'G13' = '''23'
    'D14' = '''0.001255'
    'E14' = '''-2.63%'
    'G14' = '''23'
    'D15' = '''0.04123'
    'E15' = '''-2.11%'
    'G15' = '''23'
    'D16' = '''0.005905'
    'E16' = '''-0.03%'
    'G16' = '''23'
    'B17' = 'HotbitToken'
    'C17' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'D17' = '''0.004414'
    'E17' = '''2.73%'
    'G17' = '''23'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '''3.359'
    'E18' = '''-0.32%'
    'G18' = '''23'
    'B19' = 'GateToken'
    'C19' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D19' = '''4.434'
    'E19' = '''2.40%'
    'G19' = '''23'
    'B20' = 'BTSEToken'
    'C20' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D20' = '''2.612'
    'E20' = '''2.22%'
    'G20' = '''23'
    'B21' = 'BitpandaEcosystemToken'
    'C21' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D21' = '''0.3343'
    'E21' = '''-3.80%'
    'G21' = '''23'
    'B22' = 'ProBitToken'
    'C22' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D22' = '''0.1383'
    'E22' = '''-1.84%'
    'G22' = '''23'
    'B23' = 'ZBToken'
    'C23' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D23' = '''0.2488'
    'E23' = '''-0.87%'
    'G23' = '''23'
    'B24' = 'BitKan'
    'C24' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'D24' = '''0.001302'
    'E24' = '''3.86%'
    'G24' = '''23'
    'D25' = '''0.0001281'
    'E25' = '''7.02%'
    'G25' = '''23'
    'D26' = '''0.0003737'
    'E26' = '''-0.50%'
    'G26' = '''23'
    'G27' = '''23'
    'G28' = '''23'
    'G29' = '''23'
    'G30' = '''23'
    'G31' = '''23'
    'G32' = '''23'
    'G33' = '''23'
    'G34' = '''23'
    'G35' = '''23'
    'G36' = '''23'
    'G37' = '''23'
    'D38' = '''0.02796'
    'E38' = '''9.12%'
    'G38' = '''23'
    'D39' = '''0.05772'
    'E39' = '''5.55%'
    'G39' = '''23'
    'D40' = '''0.007650'
    'E40' = '''1.15%'
    'G40' = '''23'
    'D41' = '''0.1424'
    'E41' = '''2.43%'
    'G41' = '''23'
    'D42' = '''0.007567'
    'E42' = '''1.10%'
    'G42' = '''23'
    'D43' = '''0.002001'
    'E43' = '''-1.26%'
    'G43' = '''23'
    'D44' = '''0.008033'
    'E44' = '''-3.46%'
    'G44' = '''23'
    'D45' = '''0.00007016'
    'E45' = '''-1.29%'
    'G45' = '''23'
    'D46' = '''0.00000000749'
    'E46' = '''-0.64%'
    'G46' = '''23'
    'D47' = '''0.0005796'
    'E47' = '''-0.27%'
    'G47' = '''23'
    'D48' = '''0.003533'
    'E48' = '''-26.95%'
    'G48' = '''23'
    'G49' = '''23'
    'D50' = '''0.00002098'
    'E50' = '''-0.64%'
    'G50' = '''23'
    'D51' = '''0.0001998'
    'E51' = '''-0.64%'
    'G51' = '''23'
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
